$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.386
$ws.Range("B7").Value = 5.717000000000001
$ws.Range("A9").Value = -21.33
$ws.Range("B12").Value = 5.57
$ws.Range("A13").Value = -22.116
$ws.Range("B14").Value = 6.468999999999999
$ws.Range("E15").Value = 16.417
$ws.Range("A16").Value = -21.95
$ws.Range("A18").Value = -21.93
$ws.Range("B19").Value = 8.138999999999999
$ws.Range("A20").Value = -20.257
$ws.Range("A26").Value = -21.526
$ws.Range("B26").Value = 6.018000000000001
$ws.Range("A27").Value = -21.37
$ws.Range("B27").Value = 5.694999999999999
$ws.Range("E28").Value = 17.077
$ws.Range("A29").Value = -21.752
$ws.Range("B29").Value = 6.425
$ws.Range("E33").Value = 17.288
$ws.Range("A35").Value = -20.244
$ws.Range("E35").Value = 16.522
$ws.Range("A36").Value = -21.057
$ws.Range("B37").Value = 8.402000000000001
$ws.Range("B38").Value = 5.748
$ws.Range("E38").Value = 16.531
$ws.Range("E43").Value = 17.092
$ws.Range("E44").Value = 16.559
$ws.Range("A45").Value = -21.706
$ws.Range("E45").Value = 17.009
$ws.Range("B47").Value = 5.361
$ws.Range("E47").Value = 16.621
$ws.Range("B51").Value = 6.264
$ws.Range("E51").Value = 16.934
$ws.Range("B52").Value = 5.387
$ws.Range("E54").Value = 16.64
$ws.Range("A55").Value = -21.641
$ws.Range("B55").Value = 5.531000000000001
$ws.Range("A57").Value = -22.125
$ws.Range("E57").Value = 16.331
$ws.Range("E62").Value = 16.154
$ws.Range("E63").Value = 17.67299999999999
$ws.Range("E67").Value = 17.373
$ws.Range("A69").Value = -21.584
$ws.Range("B69").Value = 5.723000000000001
$ws.Range("B70").Value = 5.140000000000001
$ws.Range("E70").Value = 17.31
$ws.Range("A76").Value = -22.125
$ws.Range("B76").Value = 5.178
$ws.Range("A78").Value = -20.242
$ws.Range("B81").Value = 5.842000000000001
$ws.Range("E81").Value = 16.638
$ws.Range("A82").Value = -22.018
$ws.Range("A83").Value = -20.237
$ws.Range("B83").Value = 7.086
$ws.Range("E88").Value = 16.363
$ws.Range("A93").Value = -21.576
$ws.Range("B94").Value = 6.465999999999999
$ws.Range("E96").Value = 16.324
$ws.Range("A97").Value = -22.136
$ws.Range("E99").Value = 16.611
$ws.Range("B100").Value = 5.199
$ws.Range("B102").Value = 6.898999999999999
